# aggiornamento a 9/09 compreso
# Appends 8 new daily rows (2021-09-02 .. 2021-09-09, Excel date serials
# 44441..44448) to the bottom of the "nuovi positivi comuni MO" table,
# extending the used range from A1:AX366 to A1:AX374.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data, one line per row: date-serial followed by the 49 municipal
# columns (B..AX), in column order.
$newRowsCsv = @"
44441,1,1,0,22,7,1,6,0,3,0,0,3,3,1,1,0,4,1,0,0,22,0,6,2,1,0,0,0,0,0,0,7,2,1,2,2,2,5,0,6,116,0,0,0,0,0,4,0,0
44442,0,1,0,5,2,4,0,3,1,0,9,4,2,0,1,0,3,1,0,9,47,2,1,0,5,1,0,0,0,0,2,11,3,0,1,0,1,6,0,3,130,1,0,0,0,0,1,0,0
44443,4,2,0,14,3,0,0,0,0,0,3,5,1,0,0,0,2,0,0,3,24,2,0,2,2,1,1,0,0,1,0,4,3,0,0,4,3,5,0,1,90,0,0,0,0,0,0,0,0
44444,0,0,0,2,6,0,0,0,0,0,3,0,1,0,0,0,2,0,0,1,26,0,0,1,3,0,0,1,1,2,1,1,4,0,0,1,1,4,0,1,62,0,0,0,0,0,0,0,0
44445,2,1,0,22,10,2,0,1,0,1,0,3,1,0,0,0,2,0,0,6,55,1,1,5,5,1,0,3,5,0,0,4,2,0,1,1,6,0,0,0,143,2,0,0,0,0,0,0,0
44446,0,2,0,5,3,0,0,2,0,0,3,0,0,0,0,0,0,0,0,5,17,0,0,0,1,0,0,2,0,2,1,5,1,1,0,1,5,2,0,3,64,0,0,0,0,0,0,0,3
44447,0,0,0,4,1,0,3,0,0,0,0,0,0,0,0,0,0,0,0,1,9,1,0,0,0,0,0,0,0,1,0,2,0,0,0,0,0,1,0,2,25,0,0,0,0,0,0,0,0
44448,0,0,1,18,2,1,5,1,0,0,1,1,1,0,1,0,0,0,2,1,18,0,0,2,15,1,0,1,0,1,1,2,2,3,0,0,0,8,0,3,94,0,0,0,0,0,0,0,2
"@

$lastDataRow = 366
$firstNewRow = 367

# The source row carries the formatting we want to replicate down the new
# rows (bold/centered/bordered date cell in column A, plain numeric cells
# elsewhere) - copy it down first so every new row inherits that style,
# then overwrite the values.
$srcRow = $ws.Range("A$lastDataRow" + ":AX$lastDataRow")

$lines = $newRowsCsv -split "`n"
$r = $firstNewRow
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }

    $dstRow = $ws.Range("A$r" + ":AX$r")
    $srcRow.Copy($dstRow)

    $values = $line -split ","
    for ($c = 1; $c -le $values.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = [double]$values[$c - 1]
    }

    $r = $r + 1
}
